$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the edited range to Text format so numeric-looking strings
# (e.g. "42.67", "61.934.55") are stored as text, matching the source data.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '61.934.55'
$ws.Range('E2').Value = '  -1.12%  '
$ws.Range('D3').Value = '3.411.39'
$ws.Range('E3').Value = '  -0.49%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '408.96'
$ws.Range('E5').Value = '  +0.50%  '
$ws.Range('D6').Value = '128.98'
$ws.Range('E6').Value = '  -0.96%  '
$ws.Range('D7').Value = '0.631'
$ws.Range('E7').Value = '  +6.08%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '0.731'
$ws.Range('E9').Value = '  +5.62%  '
$ws.Range('E10').Value = '  +1.56%  '
$ws.Range('D11').Value = '42.67'
$ws.Range('E11').Value = '  +1.77%  '
$ws.Range('D12').Value = '9.17'
$ws.Range('E12').Value = '  +9.24%  '
$ws.Range('D13').Value = '0.0000216'
$ws.Range('E13').Value = '  +38.79%  '
$ws.Range('E14').Value = '  -0.17%  '
$ws.Range('D15').Value = '3.953.80'
$ws.Range('E15').Value = '  -0.43%  '
$ws.Range('D16').Value = '21.21'
$ws.Range('E16').Value = '  +7.09%  '
$ws.Range('D17').Value = '3.410.01'
$ws.Range('E17').Value = '  +0.08%  '
$ws.Range('D18').Value = '12.47'
$ws.Range('E18').Value = '  +7.99%  '
$ws.Range('D19').Value = '1.09'
$ws.Range('E19').Value = '  +7.17%  '
$ws.Range('D20').Value = '61.931.12'
$ws.Range('E20').Value = '  -0.92%  '
$ws.Range('D21').Value = '449.27'
$ws.Range('E21').Value = '  +43.73%  '
$ws.Range('D22').Value = '91.74'
$ws.Range('E22').Value = '  +8.82%  '
$ws.Range('E23').Value = '  +1.35%  '
$ws.Range('D24').Value = '13.15'
$ws.Range('E24').Value = '  +2.78%  '
$ws.Range('E25').Value = '  +3.63%  '
$ws.Range('D26').Value = '9.36'
$ws.Range('E26').Value = '  +15.60%  '
$ws.Range('D27').Value = '33.00'
$ws.Range('E27').Value = '  +11.14%  '
$ws.Range('D29').Value = '7.62'
$ws.Range('E29').Value = '  -1.87%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '2.70'
$ws.Range('E30').Value = '  -1.33%  '
$ws.Range('B31').Value = 'Cosmos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D31').Value = '12.06'
$ws.Range('E31').Value = '  +6.08%  '
$ws.Range('D32').Value = '0.171'
$ws.Range('E32').Value = '  -1.27%  '
$ws.Range('E33').Value = '  -0.12%  '
$ws.Range('D34').Value = '42.68'
$ws.Range('E34').Value = '  -4.64%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.08%  '
$ws.Range('E36').Value = '  +3.86%  '
$ws.Range('D37').Value = '53.82'
$ws.Range('E37').Value = '  +3.89%  '
$ws.Range('E38').Value = '  -0.15%  '
$ws.Range('E39').Value = '  +2.10%  '
$ws.Range('E40').Value = '  +7.24%  '
$ws.Range('D41').Value = '2.96'
$ws.Range('E41').Value = '  -0.43%  '
$ws.Range('E42').Value = '  -2.24%  '
$ws.Range('D43').Value = '142.75'
$ws.Range('E43').Value = '  +0.41%  '
$ws.Range('D44').Value = '4.25'
$ws.Range('E44').Value = '  +8.15%  '
$ws.Range('D45').Value = '2.56'
$ws.Range('E45').Value = '  +15.59%  '
$ws.Range('D46').Value = '2.00'
$ws.Range('E46').Value = '  +1.10%  '
$ws.Range('D47').Value = '16.60'
$ws.Range('E47').Value = '  -1.11%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '0.147'
$ws.Range('E48').Value = '  +22.86%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '22.41'
$ws.Range('E49').Value = '  +5.68%  '
$ws.Range('E50').Value = '  +8.70%  '
$ws.Range('D51').Value = '3.757.96'
$ws.Range('E51').Value = '  -0.58%  '

# Remove the temporary text-format style so the cells keep their original
# (unstyled) appearance, same as the rest of the sheet.
$ws.Range("B2:E51").ClearFormats()
